# Update "Ready for handoff" status and timestamps across sheets, and
# widen the "Status" column on each sheet (Overview E:F, zh-cn C, de-de C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Overview "Latest HO Xliff Generate Date" (G2) timestamp bump ---
$overview.Range("G2").Value = "2018-02-28 11:39:50"

# --- de-de "Latest Handoff Datetime" (H2) timestamp bump ---
$dede.Range("H2").Value = "2018-02-28 11:39:50"

# --- zh-cn "Latest Handoff Datetime" (H2) timestamp bump ---
$zhcn.Range("H2").Value = "2018-02-28 11:39:28"

# --- Error Detail (R2) latest handback md commit hash bump, both langs ---
$zhcn.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/53bcfd6ddceeb0ab831545734ad52258c35fee06/e2e/8509edc2-8290-4d8e-bae0-1a55c178ac8b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/11fe1ac3f33e94da562da1ef167f46344a9f4114/e2e/8509edc2-8290-4d8e-bae0-1a55c178ac8b.md."
$dede.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/53bcfd6ddceeb0ab831545734ad52258c35fee06/e2e/8509edc2-8290-4d8e-bae0-1a55c178ac8b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/11fe1ac3f33e94da562da1ef167f46344a9f4114/e2e/8509edc2-8290-4d8e-bae0-1a55c178ac8b.md."

# --- Widen "Status" columns to fit the new longer text ---
# (17.2159881591797 characters is the target stored width; the COM layer
# quantizes ColumnWidth assignments onto a 1/6-character grid, so feed it
# an input that lands in the bucket nearest that target.)
$overview.Range("E1:F1").ColumnWidth = 16.33
$zhcn.Range("C1").ColumnWidth = 16.33
$dede.Range("C1").ColumnWidth = 16.33
